{"js": "// Remove 4 of the blank paragraphs that were reserved right after the\n// \"Z wykresu wynika...\" sentence (empty space kept for drawing \"Wykres 6\").\nconst body = context.document.body;\n\n// Locate the sentence paragraph via search rather than a hard-coded index.\nconst results = body.search(\"Z wykresu wynika\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Anchor sentence 'Z wykresu wynika...' not found\");\n}\n\nconst anchorParagraph = results.items[0].paragraphs.getFirst();\nawait context.sync();\n\n// Delete the 4 empty paragraphs immediately following the sentence.\nlet cur = anchorParagraph;\nfor (let i = 0; i < 4; i++) {\n  const next = cur.getNext();\n  next.delete();\n  await context.sync();\n}\n", "ps1": "# Remove 4 of the blank paragraphs that were reserved after the\n# \"Z wykresu wynika...\" sentence (whitespace kept for drawing \"Wykres 6\").\n$d = $word.ActiveDocument\n\n# Locate the sentence paragraph robustly via Find (diacritics-free needle\n# is enough to get a unique match) instead of a hard-coded paragraph index.\n$rng = $d.Content\n$found = $rng.Find.Execute(\"Z wykresu wynika, ze im dluzsze dlugosc\")\nif (-not $found) {\n    $rng = $d.Content\n    $found = $rng.Find.Execute(\"Z wykresu wynika\")\n}\n\n$anchorParagraph = $rng.Paragraphs(1)\n$anchorIndex = $anchorParagraph.Index\n\n# Delete the 4 empty paragraphs immediately following the sentence.\nfor ($i = 0; $i -lt 4; $i++) {\n    $target = $d.Paragraphs($anchorIndex + 1)\n    $target.Range.Delete()\n}\n"}
